$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 167, shifting existing rows 167-176 down to 168-177.
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new record (matches the
# pattern of the surrounding rows for the constant columns).
$ws.Range("A167").Value = 4
$ws.Range("B167").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C167").Value = "Los Lagos"
$ws.Range("D167").Value = 44516
$ws.Range("E167").Value = 10
$ws.Range("F167").Value = 100112043
$ws.Range("G167").Value = "Pepino ensalada"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 400
$ws.Range("K167").Value = 11000
$ws.Range("L167").Value = 11000
$ws.Range("M167").Value = 11000
$ws.Range("N167").Value = "$/caja 60 unidades"
$ws.Range("O167").Value = "Región de Arica y Parinacota"
$ws.Range("P167").Value = 183
$ws.Range("Q167").Value = 60
$ws.Range("R167").Value = "Hortaliza"
